# g3.5 - correção do nome da aba e inclusão do ano na planilha

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Correct the worksheet (tab) name
$ws.Name = "g3.5b Média"

# 2) Add the "Ano" column (D) with the same header style as the other headers
$ws.Range("A1").Copy($ws.Range("D1"))
$ws.Range("D1").Value = "Ano"

# 3) Fill the "Ano" values for the data rows (2-10)
for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 4).Value = "2010-2023"
}
